# Updates the bisection-method table on Sheet1 ("tabla_biseccion"):
#   - Recomputes the xn / fxn / E columns (B:D) for rows 2-21 using the
#     real bisection data (f(x) = x^3 - x - 2, bracket [1,2]).
#   - Removes the former row 22 (the table now only needs 20 iterations),
#     which also shrinks the used range from A1:D22 to A1:D21.
#
# Cell values are written as literal text (matching the workbook's existing
# "numbers stored as text" convention) by forcing Text format before the
# assignment, then restoring each cell's original style so no unrelated
# formatting/style metadata is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$data = @{
    2  = @("1.5",               "-0.125",               "1.000001")
    3  = @("1.75",               "1.609375",             "0.142857142857143")
    4  = @("1.625",              "0.666015625",          "0.0769230769230769")
    5  = @("1.5625",             "0.252197265625",       "0.04")
    6  = @("1.53125",            "0.059112548828125",    "0.0204081632653061")
    7  = @("1.515625",           "-0.0340538024902344",  "0.0103092783505155")
    8  = @("1.5234375",          "0.0122504234313965",   "0.0051282051282051")
    9  = @("1.51953125",         "-0.0109712481498718",  "0.0025706940874036")
    10 = @("1.521484375",        "0.0006221756339073",   "0.0012836970474967")
    11 = @("1.5205078125",       "-0.0051788864657282",  "0.0006422607578676")
    12 = @("1.52099609375",      "-0.0022794433170929",  "0.0003210272873194")
    13 = @("1.521240234375",     "-0.0008289058605441",  "0.0001604878831648")
    14 = @("1.5213623046875",    "-0.0001034331235132",  "8.02375030089064e-05")
    15 = @("1.52142333984375",   "0.0002593542519662",   "4.01171420548e-05")
    16 = @("1.52139282226562",   "7.79563135040462e-05", "2.00589733817423e-05")
    17 = @("1.52137756347656",   "-1.27394676745496e-05","1.00295872824833e-05")
    18 = @("1.52138519287109",   "3.2608157245928e-05",  "5.01476849321251e-06")
    19 = @("1.52138137817383",   "9.93427836881722e-06", "2.50739053359778e-06")
    20 = @("1.5213794708252",    "-1.40261125736174e-06","1.25369683855268e-06")
    21 = @("1.52138042449951",   "4.26582940482589e-06", "6.26848026337647e-07")
}

foreach ($r in 2..21) {
    $vals = $data[$r]

    $cellB = $ws.Cells.Item($r, 2)
    Set-TextValue $cellB $vals[0]

    $cellC = $ws.Cells.Item($r, 3)
    Set-TextValue $cellC $vals[1]

    $cellD = $ws.Cells.Item($r, 4)
    Set-TextValue $cellD $vals[2]
}

# Drop the old row 22 (its data is gone now that the table only runs to
# iteration 19); this also recomputes the sheet's dimension to A1:D21.
$ws.Rows(22).Delete()
